$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the view: scroll back to top-left (A1) and move the active selection to H4
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H4").Select()

# Update row 42 data: thiscolor -> red.png, target -> 1, corrAns -> 1
$ws.Range("B42").Value = "red.png"
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 1
